$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values changed
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2 and D2 cleared (deleted), C2 and E2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -4.5192477786255836
$ws.Range("D2").ClearContents()
$ws.Range("E2").Value = -5.1444421971330101

# Row 3: updated values
$ws.Range("B3").Value = -5.4378866419480616
$ws.Range("C3").Value = -0.88642087638426403
$ws.Range("D3").Value = -8.3587685308973292
$ws.Range("E3").Value = 9.4870072226013917

# Selection reflects the new, smaller highlighted range
$ws.Range("B1:E3").Select()
